$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows (for the new materials "Holden" and "Rizzie Spiral")
#    right after row 3 ("Spiral5"). This pushes the existing data for rows
#    4-29 down to rows 6-31, automatically carrying their values/labels and
#    growing the used range to A1:T31.
# ---------------------------------------------------------------------------
$ws.Rows("4:5").Insert()

# Re-apply the column-A number formatting (bold / bordered / centered style)
# that Insert() does not fully carry over to the freshly inserted blank rows.
$ws.Range("A2").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the row labels and simulation data for the two new materials.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$row4Vals = @(0.9681978817623736,1.004678467702034,1.013957874257259,0.9681978817623736,1.011044505675801,0.9844682076787383,1.013957874257259,1.007383834233592,1.013957874257259,1.004678467702034,0.986438174732204,0.986438174732204,0.9857815190477154,0.9956114079072225,0.9956114079072226,1.000198024494732,1.000198024494732,0.9982884618849663)
$row5Vals = @(1.319937262685466,0.9144853086353874,0.9307118939630491,1.319937262685466,0.7709887543065238,1.178534735263661,0.9307118939630491,0.9192160884878839,0.9307118939630491,0.9144853086353874,1.117211285660427,1.117211285660427,1.137652435528172,1.055044821761301,1.055044821761301,1.023961589811738,1.023961589811738,1.005645673890329)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $row4Vals[$i]
    $ws.Range($cols[$i] + "5").Value = $row5Vals[$i]
}

# ---------------------------------------------------------------------------
# 3. Rename "Thomas Hex" -> "Matthies Hex" (row shifted from 9 to 11 by the
#    insert above).
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Matthies Hex"
